$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    4 = @{
        B = 0.5989399169560075
        C = 0.604
        D = 0.6029096595738366
        E = 0.601
        F = 0.5703248346767411
        G = 0.639
        H = 0.5300924034930301
        I = 0.531
        J = 0.7076695793592165
        K = 0.7709999999999999
        L = 0.6631644263236074
        M = 0.6815
    }
    5 = @{
        B = 0.6498387506782149
        C = 0.631
        D = 0.6934547866226503
        E = 0.669
        F = 0.5007830588947616
        G = 0.4819999999999999
        H = 0.5408701106508991
        I = 0.5405
        J = 0.683940924245246
        K = 0.67
        L = 0.7252685055610366
        M = 0.7105
    }
    6 = @{
        B = 0.6645574225582199
        C = 0.68
        D = 0.6672017806586634
        E = 0.6679999999999999
        F = 0.5428036421342159
        G = 0.52
        H = 0.5908772466563478
        I = 0.5769999999999998
        J = 0.7223424315728552
        K = 0.731
        L = 0.7402352666560535
        M = 0.724
    }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
